# Fill in personal/form details for Mr. Md. Hasibul Haque (Math) on Sheet1,
# widen column A so the footer labels still fit, and move the saved
# selection/scroll position, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Name / designation / department (row 3-5 form fields) ---
$ws.Range("A3").Value = "নাম: Mr. Md. Hasibul Haque (Math)"
$ws.Range("A4").Value = "পদবী: সহকারী অধ্যাপক"
$ws.Range("G4").Value = "৪র্থ"
$ws.Range("I4").Value = "১ম"
$ws.Range("B5").Value = "সিএসই"
$ws.Range("F5").Value = "বিভাগ :গণিত"

# --- Amount in words (row 32) ---
$ws.Range("A32").Value = "কথায়:দুই হাজার সাতশো টাকা মাত্র।"

# --- Column A widened to fit the longer name/content ---
$ws.Columns.Item(1).ColumnWidth = 14.33203125

# --- Restore the saved scroll position / selection ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("I32").Select()
